$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: cell address, new text, and whether the text looks like a plain
# number (so it must be forced into Text format to avoid Excel silently
# converting it to a numeric value and dropping things like trailing zeros).
$updates = @(
    ,@("D2", "30.698.82", $false)
    ,@("D3", "2.109.12", $false)
    ,@("E3", "  +10.32%  ", $false)
    ,@("E4", "  -0.02%  ", $false)
    ,@("D5", "330.29", $true)
    ,@("E5", "  +3.09%  ", $false)
    ,@("D6", "0.9996", $true)
    ,@("E6", "  -0.04%  ", $false)
    ,@("D7", "0.5216", $true)
    ,@("E7", "  +3.01%  ", $false)
    ,@("D8", "0.4414", $true)
    ,@("E8", "  +8.23%  ", $false)
    ,@("D9", "0.09019", $true)
    ,@("E9", "  +7.98%  ", $false)
    ,@("D10", "46.49", $true)
    ,@("E10", "  +9.55%  ", $false)
    ,@("D11", "1.176", $true)
    ,@("E11", "  +6.35%  ", $false)
    ,@("D12", "24.88", $true)
    ,@("E12", "  +4.60%  ", $false)
    ,@("D13", "2.105.26", $false)
    ,@("E13", "  +10.47%  ", $false)
    ,@("D14", "6.796", $true)
    ,@("E14", "  +6.22%  ", $false)
    ,@("D15", "7.699", $true)
    ,@("E15", "  +6.44%  ", $false)
    ,@("D16", "98.07", $true)
    ,@("E16", "  +6.04%  ", $false)
    ,@("D17", "1.000", $true)
    ,@("E17", "  -0.19%  ", $false)
    ,@("E18", "  +3.85%  ", $false)
    ,@("D19", "0.06622", $true)
    ,@("E19", "  +1.80%  ", $false)
    ,@("D20", "19.25", $true)
    ,@("E20", "  +4.27%  ", $false)
    ,@("E21", "  +8.12%  ", $false)
    ,@("D22", "0.9997", $true)
    ,@("E22", "  -0.07%  ", $false)
    ,@("D23", "30.846.80", $false)
    ,@("E23", "  +2.38%  ", $false)
    ,@("D24", "12.03", $true)
    ,@("E24", "  +6.20%  ", $false)
    ,@("D25", "2.356.56", $false)
    ,@("E25", "  +10.81%  ", $false)
    ,@("D26", "2.255", $true)
    ,@("E26", "  +2.94%  ", $false)
    ,@("D27", "22.95", $true)
    ,@("E27", "  +5.36%  ", $false)
    ,@("B28", "LidoDAOToken", $false)
    ,@("C28", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", $false)
    ,@("D28", "2.543", $true)
    ,@("E28", "  +11.46%  ", $false)
    ,@("B29", "Monero", $false)
    ,@("C29", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", $false)
    ,@("D29", "163.33", $true)
    ,@("E29", "  +0.16%  ", $false)
    ,@("D30", "134.23", $true)
    ,@("E30", "  +4.28%  ", $false)
    ,@("D31", "1.189", $true)
    ,@("E31", "  +3.90%  ", $false)
    ,@("D32", "0.1070", $true)
    ,@("E32", "  +2.65%  ", $false)
    ,@("D33", "6.228", $true)
    ,@("E33", "  +4.55%  ", $false)
    ,@("B34", "ARBITRUM", $false)
    ,@("C34", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", $false)
    ,@("D34", "1.532", $true)
    ,@("E34", "  +28.13%  ", $false)
    ,@("B35", "HuobiToken", $false)
    ,@("C35", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", $false)
    ,@("D35", "3.905", $true)
    ,@("E35", "  +2.98%  ", $false)
    ,@("D36", "0.02585", $true)
    ,@("E36", "  +5.25%  ", $false)
    ,@("D37", "5.620", $true)
    ,@("E37", "  +4.88%  ", $false)
    ,@("D38", "0.06746", $true)
    ,@("E38", "  +5.56%  ", $false)
    ,@("D39", "9.523", $true)
    ,@("E39", "  +10.22%  ", $false)
    ,@("E40", "  +12.27%  ", $false)
    ,@("D41", "0.2251", $true)
    ,@("E41", "  +4.64%  ", $false)
    ,@("D42", "0.6791", $true)
    ,@("D43", "1.252", $true)
    ,@("E43", "  +3.26%  ", $false)
    ,@("D44", "14.25", $true)
    ,@("E44", "  +6.43%  ", $false)
    ,@("D45", "0.9993", $true)
    ,@("E45", "  -0.04%  ", $false)
    ,@("D46", "0.6325", $true)
    ,@("E46", "  +3.96%  ", $false)
    ,@("E47", "  +2.89%  ", $false)
    ,@("D48", "3.653", $true)
    ,@("E48", "  +0.85%  ", $false)
    ,@("D49", "1.281", $true)
    ,@("E49", "  +5.88%  ", $false)
    ,@("D50", "124.02", $true)
    ,@("E50", "  +1.48%  ", $false)
    ,@("D51", "83.05", $true)
    ,@("E51", "  +5.11%  ", $false)
)

foreach ($u in $updates) {
    $cell = $u[0]
    $value = $u[1]
    $isNumeric = $u[2]
    $range = $ws.Range($cell)
    if ($isNumeric) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
